$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($ws, $addr, $text) {
    $c = $ws.Range($addr)
    $c.Formula = "'" + $text
    $c.Style = "Normal"
}

$updates = @(
    @("D2", '69.606.79'),
    @("E2", '  -0.34%  '),
    @("D3", '3.789.61'),
    @("E3", '  +0.80%  '),
    @("E4", '  -0.06%  '),
    @("D5", '614.98'),
    @("E5", '  -0.93%  '),
    @("D6", '177.74'),
    @("E6", '  -2.22%  '),
    @("D7", '3.792.12'),
    @("E7", '  +0.91%  '),
    @("E8", '  -0.03%  '),
    @("D9", '0.526'),
    @("E9", '  -1.52%  '),
    @("E10", '  -1.37%  '),
    @("E11", '  +1.98%  '),
    @("D12", '0.485'),
    @("E12", '  -1.44%  '),
    @("D13", '39.90'),
    @("E13", '  -3.62%  '),
    @("D14", '0.0000255'),
    @("E14", '  -2.20%  '),
    @("D15", '4.412.64'),
    @("E15", '  +0.78%  '),
    @("D16", '3.784.79'),
    @("E16", '  +0.52%  '),
    @("D17", '69.666.19'),
    @("E17", '  -0.41%  '),
    @("D18", '7.57'),
    @("E18", '  -0.59%  '),
    @("E19", '  -3.55%  '),
    @("D20", '510.93'),
    @("E20", '  +0.39%  '),
    @("D21", '16.63'),
    @("E21", '  -0.79%  '),
    @("D22", '9.56'),
    @("E22", '  -0.65%  '),
    @("D23", '0.736'),
    @("E23", '  +0.75%  '),
    @("D24", '2.49'),
    @("E24", '  -1.61%  '),
    @("D25", '86.42'),
    @("E25", '  -1.16%  '),
    @("D26", '12.91'),
    @("D27", '0.0000142'),
    @("E27", '  +3.85%  '),
    @("D28", '10.58'),
    @("E28", '  -5.12%  '),
    @("E29", '  +0.00%  '),
    @("E30", '  +3.38%  '),
    @("E31", '  -0.69%  '),
    @("E32", '  +2.72%  '),
    @("D33", '31.33'),
    @("E33", '  +0.37%  '),
    @("E34", '  -0.27%  '),
    @("E35", '  -0.17%  '),
    @("E36", '  -1.36%  '),
    @("D37", '6.16'),
    @("E37", '  -1.08%  '),
    @("E38", '  +6.70%  '),
    @("D39", '477.35'),
    @("E39", '  +11.19%  '),
    @("D40", '0.341'),
    @("E40", '  +0.72%  '),
    @("E41", '  -2.51%  '),
    @("D42", '3.01'),
    @("E42", '  +5.64%  '),
    @("D43", '49.78'),
    @("E43", '  -0.87%  '),
    @("D44", '44.31'),
    @("E44", '  -3.07%  '),
    @("D45", '8.60'),
    @("E45", '  -1.93%  '),
    @("D46", '2.949.59'),
    @("E46", '  -2.03%  '),
    @("D47", '0.0363'),
    @("E47", '  -0.79%  '),
    @("D48", '27.57'),
    @("E48", '  +0.14%  '),
    @("D49", '139.66'),
    @("E49", '  +1.96%  '),
    @("E50", '  +0.05%  '),
    @("E51", '  -2.17%  ')
)

foreach ($u in $updates) {
    Set-CellText $ws $u[0] $u[1]
}
